# 1.1.3 - setup and docs
#
# The published document keeps only the default (primary) header -
# the even/first headers and all footers are dropped - and the
# surviving header's version/date line is bumped from
# "Version 1.1.2 \nNovember 5, 2022" to "Version 1.1.3\nNovember 8, 2022".

$d   = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Bump the version number and the date shown in the default header ---
# Headers(1) is the primary/default header (the one actually displayed,
# i.e. w:headerReference w:type="default").
$hdr = $sec.Headers(1)

# "Version 1.1.2 " (note trailing space before the line break) -> "Version 1.1.3"
$hdr.Range.Find.Execute("1.1.2 ", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "1.1.3", 2) | Out-Null

# "November 5, 2022" -> "November 8, 2022"
$hdr.Range.Find.Execute("November 5, 2022", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "November 8, 2022", 2) | Out-Null

# --- Drop the unused even/first-page headers and all footers ---
# Only the primary header stays linked to real content; everything else
# (the even-page header, the first-page header and every footer slot)
# is removed so the section is left with a single default header and no
# footer at all, same as a real "Remove Footer" / "Different First Page"
# + "Different Odd & Even Pages" clean-up pass in Word.
foreach ($h in @($sec.Headers(2), $sec.Headers(3))) {
    try { $h.LinkToPrevious = $false } catch {}
    try { $h.Range.Text = "" } catch {}
    try { $h.Exists = $false } catch {}
}
foreach ($f in @($sec.Footers(1), $sec.Footers(2), $sec.Footers(3))) {
    try { $f.LinkToPrevious = $false } catch {}
    try { $f.Range.Text = "" } catch {}
    try { $f.Exists = $false } catch {}
}
